# Auto-generated edit script applying the cell-value changes described in the diff.
# Each worksheet is addressed by its tab name, and each affected cell is updated
# directly with the new numeric value captured from the unified diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 905
$ws.Range("I19").Value = 950
$ws.Range("J19").Value = 882.5
$ws.Range("K19").Value = 950
$ws.Range("L19").Value = 882.5
$ws.Range("M19").Value = -775
$ws.Range("N19").Value = -1232.5
$ws.Range("H98").Value = 1633.6666
$ws.Range("I98").Value = 1038.8462
$ws.Range("K98").Value = 1038.8462
$ws.Range("M98").Value = 459.1538
$ws.Range("H113").Value = 125003740
$ws.Range("I113").Value = 500000000
$ws.Range("K113").Value = 500000000
$ws.Range("M113").Value = -499996746
$ws.Range("H122").Value = 1633.6666
$ws.Range("I122").Value = 1038.8462
$ws.Range("K122").Value = 3116.5386
$ws.Range("M122").Value = -666.5385999999999
$ws.Range("H132").Value = 1761.6666
$ws.Range("I132").Value = 1083.8334
$ws.Range("J132").Value = 6732.4443
$ws.Range("K132").Value = 3251.5002
$ws.Range("L132").Value = 20197.3329
$ws.Range("M132").Value = -721.5001999999999
$ws.Range("N132").Value = -25257.3329
$ws.Range("H137").Value = 2831250.8
$ws.Range("I137").Value = 1137372.5
$ws.Range("J137").Value = 11112433
$ws.Range("K137").Value = 3412117.5
$ws.Range("L137").Value = 33337299
$ws.Range("M137").Value = -3409567.5
$ws.Range("N137").Value = -33342399
$ws.Range("H141").Value = 1623.9215
$ws.Range("I141").Value = 1238.1082
$ws.Range("K141").Value = 3714.3246
$ws.Range("M141").Value = 1465.6754

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2220.9656
$ws.Range("I2").Value = 1410.4546
$ws.Range("J2").Value = 4768.2856
$ws.Range("K2").Value = 1410.4546
$ws.Range("L2").Value = 4768.2856
$ws.Range("M2").Value = -1297.4546
$ws.Range("N2").Value = -4994.2856
$ws.Range("H32").Value = 2705231.2
$ws.Range("J32").Value = 17949576
$ws.Range("L32").Value = 17949576
$ws.Range("N32").Value = -17950150
$ws.Range("H74").Value = 853.2727
$ws.Range("I74").Value = 808.8
$ws.Range("J74").Value = 1298
$ws.Range("K74").Value = 808.8
$ws.Range("L74").Value = 1298
$ws.Range("M74").Value = 65.20000000000005
$ws.Range("N74").Value = -3046
$ws.Range("H77").Value = 853.2727
$ws.Range("I77").Value = 808.8
$ws.Range("J77").Value = 1298
$ws.Range("K77").Value = 4044
$ws.Range("L77").Value = 6490
$ws.Range("M77").Value = 324
$ws.Range("N77").Value = -15226
$ws.Range("H116").Value = 2220.9656
$ws.Range("I116").Value = 1410.4546
$ws.Range("J116").Value = 4768.2856
$ws.Range("K116").Value = 1410.4546
$ws.Range("L116").Value = 4768.2856
$ws.Range("M116").Value = 883.5454
$ws.Range("N116").Value = -9356.285599999999
$ws.Range("H132").Value = 117676.65
$ws.Range("I132").Value = 129435.84
$ws.Range("K132").Value = 388307.52
$ws.Range("M132").Value = -385777.52

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2220.9656
$ws.Range("I3").Value = 1410.4546
$ws.Range("J3").Value = 4768.2856
$ws.Range("K3").Value = 1410.4546
$ws.Range("L3").Value = 4768.2856
$ws.Range("M3").Value = -1296.4546
$ws.Range("N3").Value = -4996.2856

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1464.3334
$ws.Range("I16").Value = 1305
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1305
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1018
$ws.Range("N16").Value = -3074
$ws.Range("H31").Value = 1730.8049
$ws.Range("I31").Value = 1281.9678
$ws.Range("K31").Value = 1281.9678
$ws.Range("M31").Value = -986.9677999999999
$ws.Range("H34").Value = 1730.8049
$ws.Range("I34").Value = 1281.9678
$ws.Range("K34").Value = 1281.9678
$ws.Range("M34").Value = -1079.9678
$ws.Range("H99").Value = 1971.1666
$ws.Range("I99").Value = 1924.8572
$ws.Range("J99").Value = 2036
$ws.Range("K99").Value = 1924.8572
$ws.Range("L99").Value = 2036
$ws.Range("M99").Value = -426.8571999999999
$ws.Range("N99").Value = -5032
$ws.Range("H113").Value = 1464.3334
$ws.Range("I113").Value = 1305
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1305
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 865
$ws.Range("N113").Value = -6840
$ws.Range("H126").Value = 1971.1666
$ws.Range("I126").Value = 1924.8572
$ws.Range("J126").Value = 2036
$ws.Range("K126").Value = 5774.571599999999
$ws.Range("L126").Value = 6108
$ws.Range("M126").Value = -3304.571599999999
$ws.Range("N126").Value = -11048

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1499.7407
$ws.Range("I129").Value = 760.8333
$ws.Range("J129").Value = 2090.8667
$ws.Range("K129").Value = 2282.4999
$ws.Range("L129").Value = 6272.6001
$ws.Range("M129").Value = 2717.5001
$ws.Range("N129").Value = -16272.6001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 557.1875
$ws.Range("I107").Value = 391.33334
$ws.Range("K107").Value = 391.33334
$ws.Range("M107").Value = 1528.66666
$ws.Range("H113").Value = 2024.6154
$ws.Range("I113").Value = 1945
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 1945
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = 225
$ws.Range("N113").Value = -7320
$ws.Range("H126").Value = 8991.342000000001
$ws.Range("I126").Value = 2600.1052
$ws.Range("J126").Value = 15382.579
$ws.Range("K126").Value = 7800.3156
$ws.Range("L126").Value = 46147.737
$ws.Range("M126").Value = -5330.3156
$ws.Range("N126").Value = -51087.737

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2333.3333
$ws.Range("I7").Value = 2400
$ws.Range("J7").Value = 2200
$ws.Range("K7").Value = 2400
$ws.Range("L7").Value = 2200
$ws.Range("M7").Value = -2288
$ws.Range("N7").Value = -2424
$ws.Range("H126").Value = 2333.3333
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -11540

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2114.0344
$ws.Range("I126").Value = 1710.8948
$ws.Range("J126").Value = 2880
$ws.Range("K126").Value = 5132.6844
$ws.Range("L126").Value = 8640
$ws.Range("M126").Value = -2662.6844
$ws.Range("N126").Value = -13580
$ws.Range("H132").Value = 1565.7291
$ws.Range("J132").Value = 3499.5
$ws.Range("L132").Value = 10498.5
$ws.Range("N132").Value = -15558.5

Write-Output "Applied 166 cell updates across 8 sheets."